$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 147.17647
$ws.Range("I55").Value = 213.71428
$ws.Range("K55").Value = 213.71428
$ws.Range("M55").Value = 0.2857199999999978
$ws.Range("H70").Value = 1730
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 2075
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 6225
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -6765
$ws.Range("H73").Value = 1730
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 2075
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 6225
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -8097
$ws.Range("H107").Value = 1116.1428
$ws.Range("I107").Value = 1112.1875
$ws.Range("K107").Value = 1112.1875
$ws.Range("M107").Value = 807.8125
$ws.Range("H116").Value = 41365.965
$ws.Range("I116").Value = 80960.38
$ws.Range("K116").Value = 80960.38
$ws.Range("M116").Value = -77518.38
$ws.Range("H132").Value = 43559.68
$ws.Range("I132").Value = 47010.086
$ws.Range("K132").Value = 141030.258
$ws.Range("M132").Value = -138500.258
$ws.Range("H137").Value = 2696.04
$ws.Range("I137").Value = 2514.7144
$ws.Range("K137").Value = 7544.1432
$ws.Range("M137").Value = -4994.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3685.8235
$ws.Range("I20").Value = 3432.4285
$ws.Range("K20").Value = 3432.4285
$ws.Range("M20").Value = -3185.4285
$ws.Range("H23").Value = 20013
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H24").Value = 25059
$ws.Range("I24").Value = 2200
$ws.Range("J24").Value = 32678.666
$ws.Range("K24").Value = 2200
$ws.Range("L24").Value = 32678.666
$ws.Range("M24").Value = -1965
$ws.Range("N24").Value = -33148.666
$ws.Range("H29").Value = 20684
$ws.Range("J29").Value = 30518
$ws.Range("L29").Value = 30518
$ws.Range("N29").Value = -31096
$ws.Range("H32").Value = 10882.25
$ws.Range("J32").Value = 10882.25
$ws.Range("L32").Value = 10882.25
$ws.Range("N32").Value = -11650.25
$ws.Range("H39").Value = 21512
$ws.Range("J39").Value = 21512
$ws.Range("L39").Value = 21512
$ws.Range("N39").Value = -22290
$ws.Range("H105").Value = 2086.9473
$ws.Range("I105").Value = 2156.0588
$ws.Range("J105").Value = 1499.5
$ws.Range("K105").Value = 2156.0588
$ws.Range("L105").Value = 1499.5
$ws.Range("M105").Value = -409.0587999999998
$ws.Range("N105").Value = -4993.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2025.7727
$ws.Range("I58").Value = 2040.1666
$ws.Range("K58").Value = 2040.1666
$ws.Range("M58").Value = -1837.1666
$ws.Range("H62").Value = 10849.588
$ws.Range("I62").Value = 7332.778
$ws.Range("J62").Value = 14806
$ws.Range("K62").Value = 7332.778
$ws.Range("L62").Value = 14806
$ws.Range("M62").Value = -6708.778
$ws.Range("N62").Value = -16054
$ws.Range("H65").Value = 10849.588
$ws.Range("I65").Value = 7332.778
$ws.Range("J65").Value = 14806
$ws.Range("K65").Value = 36663.89
$ws.Range("L65").Value = 74030
$ws.Range("M65").Value = -33543.89
$ws.Range("N65").Value = -80270
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H105").Value = 2577.55
$ws.Range("I105").Value = 2108.8823
$ws.Range("K105").Value = 2108.8823
$ws.Range("M105").Value = -361.8823000000002
$ws.Range("H121").Value = 13925.143
$ws.Range("J121").Value = 13925.143
$ws.Range("L121").Value = 13925.143
$ws.Range("N121").Value = -16545.143
$ws.Range("H132").Value = 6251.8887
$ws.Range("I132").Value = 7966.3335
$ws.Range("K132").Value = 23899.0005
$ws.Range("M132").Value = -21369.0005
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059
$ws.Range("H134").Value = 2471.1333
$ws.Range("I134").Value = 2472.3333
$ws.Range("K134").Value = 7416.999899999999
$ws.Range("M134").Value = -4881.999899999999
$ws.Range("H136").Value = 2025.7727
$ws.Range("I136").Value = 2040.1666
$ws.Range("K136").Value = 6120.4998
$ws.Range("M136").Value = -3570.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 535.34375
$ws.Range("I2").Value = 644.6818
$ws.Range("J2").Value = 294.8
$ws.Range("K2").Value = 3868.0908
$ws.Range("L2").Value = 1768.8
$ws.Range("M2").Value = -3755.0908
$ws.Range("N2").Value = -1994.8
$ws.Range("H11").Value = 1417.5
$ws.Range("I11").Value = 732.3333
$ws.Range("J11").Value = 2102.6667
$ws.Range("K11").Value = 2196.9999
$ws.Range("L11").Value = 6308.000100000001
$ws.Range("M11").Value = -2056.9999
$ws.Range("N11").Value = -6588.000100000001
$ws.Range("H12").Value = 405
$ws.Range("J12").Value = 372
$ws.Range("L12").Value = 1116
$ws.Range("N12").Value = -1462
$ws.Range("H33").Value = 223.38461
$ws.Range("I33").Value = 93
$ws.Range("K33").Value = 558
$ws.Range("M33").Value = -275
$ws.Range("H39").Value = 5935.778
$ws.Range("I39").Value = 2831
$ws.Range("J39").Value = 6556.7334
$ws.Range("K39").Value = 8493
$ws.Range("L39").Value = 19670.2002
$ws.Range("M39").Value = -8199
$ws.Range("N39").Value = -20258.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5865
$ws.Range("I68").Value = 4200.9565
$ws.Range("J68").Value = 25001.5
$ws.Range("K68").Value = 4200.9565
$ws.Range("L68").Value = 25001.5
$ws.Range("M68").Value = -3451.9565
$ws.Range("N68").Value = -26499.5
$ws.Range("H71").Value = 5865
$ws.Range("I71").Value = 4200.9565
$ws.Range("J71").Value = 25001.5
$ws.Range("K71").Value = 21004.7825
$ws.Range("L71").Value = 125007.5
$ws.Range("M71").Value = -17260.7825
$ws.Range("N71").Value = -132495.5
$ws.Range("H82").Value = 3335.9048
$ws.Range("I82").Value = 2767.5293
$ws.Range("K82").Value = 2767.5293
$ws.Range("M82").Value = -2406.5293
$ws.Range("H85").Value = 3335.9048
$ws.Range("I85").Value = 2767.5293
$ws.Range("K85").Value = 2767.5293
$ws.Range("M85").Value = -1519.5293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2403.8518
$ws.Range("I113").Value = 740.5833
$ws.Range("J113").Value = 3734.4666
$ws.Range("K113").Value = 2221.7499
$ws.Range("L113").Value = 11203.3998
$ws.Range("M113").Value = -51.7498999999998
$ws.Range("N113").Value = -15543.3998
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 2431.75
$ws.Range("I136").Value = 2431.75
$ws.Range("K136").Value = 7295.25
$ws.Range("M136").Value = -4745.25
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
